# Applies the "Sync attendance_reports, modules_schedules, and assets from
# main repo - 2026-01-08 11:16:15" update to the Session Analysis Results
# sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) "Recorded By" column (G): every "<email>, System" becomes
#    "System, <email>" (order of the two recorders swapped).
# ---------------------------------------------------------------------
$used = $ws.UsedRange
$lastRow = $used.Rows.Count
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
}

# ---------------------------------------------------------------------
# 2) Class Statistics block (B1A1 group, rows 6-10): plain numbers.
# ---------------------------------------------------------------------
$ws.Range("L6").Value = 219
$ws.Range("L7").Value = 27

# ---------------------------------------------------------------------
# 3) Class Statistics percentages (stored as literal text, e.g. "68.9%")
#    plus Group Statistics percentages for B1D1..B1F2 (rows 21-26).
#    These must stay TEXT (not be re-interpreted as numeric percentages)
#    and keep their original style index, so: force a "@" (text) number
#    format before writing, then restore the cell's original formatting
#    via a Copy/PasteSpecial(formats) from an untouched reference cell
#    that already carries that exact style.
# ---------------------------------------------------------------------
$pctCells  = @("L9","L10","O21","P21","R21","S21","O22","P22","R22","S22","O23","P23","R23","S23","O24","P24","R24","S24","O25","P25","R25","S25","O26","P26","R26","S26")
$pctValues = @("68.9%","75.9%","18","3","66.7%","78.0%","18","3","66.7%","77.6%","18","3","66.7%","78.4%","17","4","63.0%","72.7%","18","3","66.7%","71.4%","18","3","66.7%","63.8%")

# O/P columns are plain numbers -> write directly.
$ws.Range("O21").Value = 18
$ws.Range("P21").Value = 3
$ws.Range("O22").Value = 18
$ws.Range("P22").Value = 3
$ws.Range("O23").Value = 18
$ws.Range("P23").Value = 3
$ws.Range("O24").Value = 17
$ws.Range("P24").Value = 4
$ws.Range("O25").Value = 18
$ws.Range("P25").Value = 3
$ws.Range("O26").Value = 18
$ws.Range("P26").Value = 3

# R/S (and L9/L10) columns hold text percentages.
$textPctCells  = @("L9","L10","R21","S21","R22","S22","R23","S23","R24","S24","R25","S25","R26","S26")
$textPctValues = @("68.9%","75.9%","66.7%","78.0%","66.7%","77.6%","66.7%","78.4%","63.0%","72.7%","66.7%","71.4%","66.7%","63.8%")

for ($i = 0; $i -lt $textPctCells.Count; $i++) {
    $rng = $ws.Range($textPctCells[$i])
    $rng.NumberFormat = "@"
    $rng.Value = $textPctValues[$i]
}

# Restore each cell's original (non-"@") style - L5 carries the same
# "Value" style (s=5) that L9/L10/O.../R.../S... originally used.
$styleRef = $ws.Range("L5")
foreach ($c in $textPctCells) {
    $styleRef.Copy()
    $ws.Range($c).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4) Six "Not Recorded" session rows flip to "Recorded" (B1D1, B1D2,
#    B1E1, B1E2, B1F1, B1F2 - session 21, date 08/01/2026): the
#    "Recorded By" / "Students" / "Status" cells get real values and the
#    row's pink "Not Recorded" fill becomes the green "Recorded" fill.
# ---------------------------------------------------------------------
$recordedRows = @(178, 205, 232, 259, 286, 313)
$hValues = @{178="17/23"; 205="26/30"; 232="14/26"; 259="21/28"; 286="21/26"; 313="18/29"}

# A row that already carries the green "Recorded" styling, used as the
# formatting donor for A:I on each flipped row.
$recordedStyleRef = $ws.Range("A2:I2")

foreach ($r in $recordedRows) {
    $ws.Range("G$r").Value = "dnasr281@gmail.com"
    $ws.Range("H$r").Value = $hValues[$r]
    $ws.Range("I$r").Value = "Recorded"

    $recordedStyleRef.Copy()
    $ws.Range("A$r" + ":I$r").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false
